# Update the "dSF" column (column F) values for the nola_aaron.xlsx
# save-data sheet. These values were recalculated after a data repull
# ("repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 2
    4  = 6
    5  = 7
    6  = 1
    7  = 2
    8  = 1
    9  = 3
    10 = -3
    11 = 2
    12 = 4
    13 = 1
    15 = 1
    16 = -2
    18 = 2
    19 = -3
    20 = 3
    21 = 4
    23 = 1
    24 = -3
    25 = 2
    26 = -3
    27 = -1
    28 = 8
    29 = -4
    30 = 1
    31 = 2
    32 = 2
    33 = -1
    35 = -3
    36 = -1
    38 = -4
    39 = 6
    40 = -3
    41 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
